$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 49
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()

# ALC row 53
$ws.Range("H53").Value = 321.07144
$ws.Range("J53").Value = 227.33333
$ws.Range("L53").Value = 227.33333
$ws.Range("N53").Value = -1501.33333

# ALC row 137
$ws.Range("H137").Value = 18571.428
$ws.Range("J137").Value = 18333.334
$ws.Range("L137").Value = 55000.00199999999
$ws.Range("N137").Value = -60100.00199999999

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2
$ws.Range("H2").Value = 2566.4443
$ws.Range("I2").Value = 1979.7333
$ws.Range("K2").Value = 1979.7333
$ws.Range("M2").Value = -1866.7333

# ARM row 32
$ws.Range("H32").Value = 26564.695
$ws.Range("I32").Value = 17770.076
$ws.Range("K32").Value = 17770.076
$ws.Range("M32").Value = -17483.076

# ARM row 97
$ws.Range("H97").Value = 892.5714
$ws.Range("I97").Value = 915.8333
$ws.Range("K97").Value = 915.8333
$ws.Range("M97").Value = -419.8333

# ARM row 116
$ws.Range("H116").Value = 2566.4443
$ws.Range("I116").Value = 1979.7333
$ws.Range("K116").Value = 1979.7333
$ws.Range("M116").Value = 314.2666999999999

# ARM row 132
$ws.Range("H132").Value = 3919.2856
$ws.Range("I132").Value = 2257.9092
$ws.Range("J132").Value = 10011
$ws.Range("K132").Value = 6773.7276
$ws.Range("L132").Value = 30033
$ws.Range("M132").Value = -4243.7276
$ws.Range("N132").Value = -35093

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3
$ws.Range("H3").Value = 2566.4443
$ws.Range("I3").Value = 1979.7333
$ws.Range("K3").Value = 1979.7333
$ws.Range("M3").Value = -1865.7333

# BSM row 60
$ws.Range("H60").Value = 67420
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 67420
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 67420
$ws.Range("N60").Value = -68618
$ws.Range("M60").ClearContents()

# BSM row 105
$ws.Range("H105").Value = 4005.5925
$ws.Range("I105").Value = 3508.2942
$ws.Range("K105").Value = 3508.2942
$ws.Range("M105").Value = -1761.2942

# BSM row 107
$ws.Range("H107").Value = 5133.9443
$ws.Range("I107").Value = 3815.182
$ws.Range("J107").Value = 7206.2856
$ws.Range("K107").Value = 3815.182
$ws.Range("L107").Value = 7206.2856
$ws.Range("M107").Value = -1895.182
$ws.Range("N107").Value = -11046.2856

$ws = $wb.Worksheets.Item("CRP")
# CRP row 15
$ws.Range("H15").Value = 560.28
$ws.Range("J15").Value = 517.34784
$ws.Range("L15").Value = 517.34784
$ws.Range("N15").Value = -857.34784

# CRP row 16
$ws.Range("H16").Value = 760.1667
$ws.Range("J16").Value = 975
$ws.Range("L16").Value = 975
$ws.Range("N16").Value = -1549

# CRP row 22
$ws.Range("H22").Value = 450
$ws.Range("I22").Value = 450
$ws.Range("K22").Value = 450
$ws.Range("M22").Value = -100

# CRP row 31
$ws.Range("H31").Value = 4276.0303
$ws.Range("I31").Value = 2629.889
$ws.Range("K31").Value = 2629.889
$ws.Range("M31").Value = -2334.889

# CRP row 34
$ws.Range("H34").Value = 4276.0303
$ws.Range("I34").Value = 2629.889
$ws.Range("K34").Value = 2629.889
$ws.Range("M34").Value = -2427.889

# CRP row 107
$ws.Range("H107").Value = 827.5714
$ws.Range("I107").Value = 264.2857
$ws.Range("J107").Value = 1109.2142
$ws.Range("K107").Value = 264.2857
$ws.Range("L107").Value = 1109.2142
$ws.Range("M107").Value = 1655.7143
$ws.Range("N107").Value = -4949.2142

# CRP row 113
$ws.Range("H113").Value = 760.1667
$ws.Range("J113").Value = 975
$ws.Range("L113").Value = 975
$ws.Range("N113").Value = -5315

# CRP row 134
$ws.Range("H134").Value = 3268.65
$ws.Range("I134").Value = 2568.75
$ws.Range("J134").Value = 4318.5
$ws.Range("K134").Value = 7706.25
$ws.Range("L134").Value = 12955.5
$ws.Range("M134").Value = -5171.25
$ws.Range("N134").Value = -18025.5

$ws = $wb.Worksheets.Item("CUL")
# CUL row 14
$ws.Range("H14").Value = 1679.6666
$ws.Range("I14").Value = 1679.6666
$ws.Range("K14").Value = 5038.9998
$ws.Range("M14").Value = -4865.9998

# CUL row 46
$ws.Range("H46").Value = 834349.2
$ws.Range("I46").Value = 890.2
$ws.Range("K46").Value = 2670.6
$ws.Range("M46").Value = -2579.6

# CUL row 55
$ws.Range("H55").Value = 78429.16
$ws.Range("J55").Value = 2012.5
$ws.Range("L55").Value = 6037.5
$ws.Range("N55").Value = -6391.5

# CUL row 97
$ws.Range("H97").Value = 74.14286
$ws.Range("J97").Value = 86
$ws.Range("L97").Value = 258
$ws.Range("N97").Value = -1250

$ws = $wb.Worksheets.Item("GSM")
# GSM row 80
$ws.Range("H80").Value = 8579.4
$ws.Range("J80").Value = 8422.111000000001
$ws.Range("L80").Value = 8422.111000000001
$ws.Range("N80").Value = -10418.111

# GSM row 83
$ws.Range("H83").Value = 8579.4
$ws.Range("J83").Value = 8422.111000000001
$ws.Range("L83").Value = 42110.55500000001
$ws.Range("N83").Value = -52094.55500000001

# GSM row 97
$ws.Range("H97").Value = 1814.963
$ws.Range("I97").Value = 1863.9584
$ws.Range("J97").Value = 1423
$ws.Range("K97").Value = 1863.9584
$ws.Range("L97").Value = 1423
$ws.Range("M97").Value = -1367.9584
$ws.Range("N97").Value = -2415

# GSM row 113
$ws.Range("H113").Value = 4124.3076
$ws.Range("I113").Value = 2779.25
$ws.Range("K113").Value = 2779.25
$ws.Range("M113").Value = -609.25

$ws = $wb.Worksheets.Item("LTW")
# LTW row 16
$ws.Range("H16").Value = 2004.8572
$ws.Range("I16").Value = 2105.6667
$ws.Range("K16").Value = 2105.6667
$ws.Range("M16").Value = -1935.6667

# LTW row 20
$ws.Range("H20").Value = 14006
$ws.Range("J20").Value = 14006
$ws.Range("L20").Value = 14006
$ws.Range("N20").Value = -14458

# LTW row 68
$ws.Range("H68").Value = 4333.3335
$ws.Range("J68").Value = 4333.3335
$ws.Range("L68").Value = 4333.3335
$ws.Range("N68").Value = -5831.3335

# LTW row 71
$ws.Range("H71").Value = 4333.3335
$ws.Range("J71").Value = 4333.3335
$ws.Range("L71").Value = 21666.6675
$ws.Range("N71").Value = -29154.6675

# LTW row 132
$ws.Range("H132").Value = 3580.0527
$ws.Range("I132").Value = 1876.75
$ws.Range("K132").Value = 5630.25
$ws.Range("M132").Value = -3100.25

# LTW row 136
$ws.Range("H136").Value = 3300
$ws.Range("I136").Value = 3300
$ws.Range("K136").Value = 9900
$ws.Range("M136").Value = -7350

$ws = $wb.Worksheets.Item("WVR")
# WVR row 101
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

# WVR row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# WVR row 136
$ws.Range("H136").Value = 49678.81
$ws.Range("I136").Value = 943.5625
$ws.Range("K136").Value = 2830.6875
$ws.Range("M136").Value = -280.6875
